$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: day 43975 (2020-05-24), 17:11 -> 23:59 ("Server Install") ---
$ws.Range("A6").Value = 43975
$ws.Range("B6").Value = 0.71597222222222223
$ws.Range("C6").Value = 0.99930555555555556
$ws.Range("E6").Value = "Server Install"
$ws.Range("E6").WrapText = $true

# --- Row 7: day 43976 (2020-05-25), 0:00 -> 0:16 ("Server Install") ---
$ws.Range("A7").Value = 43976
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0.011111111111111112
$ws.Range("E7").Value = "Server Install"
$ws.Range("E7").WrapText = $true

# --- Row 8: day 43976 (2020-05-25), 15:49 -> 23:11 ---
$ws.Range("A8").Value = 43976
$ws.Range("B8").Value = 0.65902777777777777
$ws.Range("C8").Value = 0.96597222222222223

# --- Total cell D26 now displays elapsed hours ([h]:mm:ss) ---
$ws.Range("D26").NumberFormat = "[h]:mm:ss"

# --- Printer / page setup (portrait) ---
$ws.PageSetup.Orientation = 1

# --- Move selection down to the recalculated total row ---
$ws.Range("D26").Select()
